# V3 and enhenced version
#
# Changes applied:
#  - Table 1 ("Recette" acceptance row): the three label cells
#    ("Acceptée, Sans Réserve :" / "Acceptée, Avec Réserve :" / "Refusée :")
#    go from right-aligned to center-aligned, and their adjoining
#    checkbox cells gain center alignment too, while losing the
#    trailing whitespace run that used to follow the checkbox field.
#  - Table 3 ("Conforme / Non conforme / N/A / Non exécuté" row): every
#    checkbox cell gains center alignment and loses its trailing
#    whitespace run.
#
# NOTE: cached Table/Row/Cell references go stale as soon as an earlier
# Find/Replace mutates the document (their Range offsets don't shift),
# so every cell is re-fetched fresh from $d right before it is touched.

$d = $word.ActiveDocument

$wdAlignParagraphCenter = 1

# ---------------------------------------------------------------------
# Table 1, row 7: Acceptée/Refusée acceptance line
# ---------------------------------------------------------------------

# Cell 1: "Acceptée, Sans Réserve :" label -> center (was right)
$d.Tables.Item(1).Rows.Item(7).Cells.Item(1).Range.Paragraphs.Item(1).Format.Alignment = $wdAlignParagraphCenter

# Cell 2: checkbox -> add center alignment, drop trailing "  " run
$d.Tables.Item(1).Rows.Item(7).Cells.Item(2).Range.Paragraphs.Item(1).Format.Alignment = $wdAlignParagraphCenter
$d.Tables.Item(1).Rows.Item(7).Cells.Item(2).Range.Find.Execute("  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

# Cell 3: "Acceptée, Avec Réserve :" label -> center (was right)
$d.Tables.Item(1).Rows.Item(7).Cells.Item(3).Range.Paragraphs.Item(1).Format.Alignment = $wdAlignParagraphCenter

# Cell 4: checkbox -> add center alignment, drop trailing "  " run
$d.Tables.Item(1).Rows.Item(7).Cells.Item(4).Range.Paragraphs.Item(1).Format.Alignment = $wdAlignParagraphCenter
$d.Tables.Item(1).Rows.Item(7).Cells.Item(4).Range.Find.Execute("  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

# Cell 5: "Refusée :" label -> center (was right)
$d.Tables.Item(1).Rows.Item(7).Cells.Item(5).Range.Paragraphs.Item(1).Format.Alignment = $wdAlignParagraphCenter

# Cell 6: checkbox -> add center alignment, drop trailing "  " run
$d.Tables.Item(1).Rows.Item(7).Cells.Item(6).Range.Paragraphs.Item(1).Format.Alignment = $wdAlignParagraphCenter
$d.Tables.Item(1).Rows.Item(7).Cells.Item(6).Range.Find.Execute("  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

# ---------------------------------------------------------------------
# Table 3, row 2: Conforme / Non conforme / N/A / Non exécuté checkboxes
# ---------------------------------------------------------------------

# Cell 1
$d.Tables.Item(3).Rows.Item(2).Cells.Item(1).Range.Paragraphs.Item(1).Format.Alignment = $wdAlignParagraphCenter
$d.Tables.Item(3).Rows.Item(2).Cells.Item(1).Range.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

# Cell 2
$d.Tables.Item(3).Rows.Item(2).Cells.Item(2).Range.Paragraphs.Item(1).Format.Alignment = $wdAlignParagraphCenter
$d.Tables.Item(3).Rows.Item(2).Cells.Item(2).Range.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

# Cell 3
$d.Tables.Item(3).Rows.Item(2).Cells.Item(3).Range.Paragraphs.Item(1).Format.Alignment = $wdAlignParagraphCenter
$d.Tables.Item(3).Rows.Item(2).Cells.Item(3).Range.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

# Cell 4
$d.Tables.Item(3).Rows.Item(2).Cells.Item(4).Range.Paragraphs.Item(1).Format.Alignment = $wdAlignParagraphCenter
$d.Tables.Item(3).Rows.Item(2).Cells.Item(4).Range.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

Write-Output "done"
